$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.409.21"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.804.10"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "2.059.60"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "1.819.10"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.629"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "34.365.76"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "0.0₃0773"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E24").Value = "  +5.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.85%  "
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0515"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("D35").Value = "1.361.72"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.651"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.96%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.92%  "
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "80.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.937"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "1.962.60"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").Value = "0.0₆0123"
$ws.Range("E51").Value = "  -4.85%  "
